$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 corresponds to ICSA_thou / Initial Jobless Claims - auto-updated data values
$ws.Range("E9").Value = 199000
$ws.Range("G9").Value = 364607.2796934866
$ws.Range("H9").Value = -20000
$ws.Range("I9").Value = -0.091324200913242
